$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Casa"
$ws.Range("O2").Value = "Leilão SFI Caixa"
$ws.Range("Q2").Value = 8083
$ws.Range("R2").Value = 37178

$ws.Range("A3").Value = "Casa"
$ws.Range("O3").Value = "Leilão SFI Caixa"
$ws.Range("Q3").Value = 9288
$ws.Range("R3").Value = 0

$ws.Range("A4").Value = "Casa"
$ws.Range("O4").Value = "Leilão SFI Caixa"
$ws.Range("Q4").Value = 9886
$ws.Range("R4").Value = 0

$ws.Range("A5").Value = "Casa"
$ws.Range("O5").Value = "Leilão SFI Caixa"
$ws.Range("Q5").Value = 9000
$ws.Range("R5").Value = 0

$ws.Range("A6").Value = "Casa"
$ws.Range("O6").Value = "Leilão SFI Caixa"
$ws.Range("Q6").Value = 9010
$ws.Range("R6").Value = 0

$ws.Range("A7").Value = "Apartamento"
$ws.Range("O7").Value = "Leilão SFI Caixa"
$ws.Range("Q7").Value = 17516
$ws.Range("R7").Value = 0

$ws.Range("A8").Value = "Casa"
$ws.Range("O8").Value = "Leilão SFI Caixa"
$ws.Range("Q8").Value = 10329
$ws.Range("R8").Value = 38155

$ws.Range("A9").Value = "Apartamento"
$ws.Range("O9").Value = "Leilão SFI Caixa"
$ws.Range("Q9").Value = 10292
$ws.Range("R9").Value = 40372

$ws.Range("A10").Value = "Apartamento"
$ws.Range("O10").Value = "Leilão SFI Caixa"
$ws.Range("Q10").Value = 53395
$ws.Range("R10").Value = 51204

$ws.Range("A11").Value = "Apartamento"
$ws.Range("O11").Value = "Leilão SFI Caixa"
$ws.Range("Q11").Value = 17494
$ws.Range("R11").Value = 0

$ws.Range("A12").Value = "Apartamento"
$ws.Range("O12").Value = "Leilão SFI Caixa"
$ws.Range("Q12").Value = 21556
$ws.Range("R12").Value = 0

$ws.Range("A13").Value = "Apartamento"
$ws.Range("O13").Value = "Leilão SFI Caixa"
$ws.Range("Q13").Value = 21729
$ws.Range("R13").Value = 0

$ws.Range("A14").Value = "Apartamento"
$ws.Range("O14").Value = "Leilão SFI Caixa"
$ws.Range("Q14").Value = 21682
$ws.Range("R14").Value = 0

$ws.Range("A15").Value = "Apartamento"
$ws.Range("O15").Value = "Leilão SFI Caixa"
$ws.Range("Q15").Value = 21470
$ws.Range("R15").Value = 0

$ws.Range("A16").Value = "Apartamento"
$ws.Range("O16").Value = "Leilão SFI Caixa"
$ws.Range("Q16").Value = 17623
$ws.Range("R16").Value = 0

$ws.Range("A17").Value = "Apartamento"
$ws.Range("O17").Value = "Leilão SFI Caixa"
$ws.Range("Q17").Value = 21658
$ws.Range("R17").Value = 0

$ws.Range("A18").Value = "Apartamento"
$ws.Range("O18").Value = "Leilão SFI Caixa"
$ws.Range("Q18").Value = 21493
$ws.Range("R18").Value = 0

$ws.Range("A19").Value = "Apartamento"
$ws.Range("O19").Value = "Leilão SFI Caixa"
$ws.Range("Q19").Value = 21509
$ws.Range("R19").Value = 0

$ws.Range("A20").Value = "Apartamento"
$ws.Range("O20").Value = "Leilão SFI Caixa"
$ws.Range("Q20").Value = 16772
$ws.Range("R20").Value = 0

$ws.Range("A21").Value = "Casa"
$ws.Range("O21").Value = "Leilão SFI Caixa"
$ws.Range("Q21").Value = 8982
$ws.Range("R21").Value = 0

$ws.Range("A22").Value = "Casa"
$ws.Range("O22").Value = "Leilão SFI Caixa"
$ws.Range("Q22").Value = 8567
$ws.Range("R22").Value = 37002

$ws.Range("A23").Value = "Casa"
$ws.Range("O23").Value = "Venda Direta Caixa"
$ws.Range("Q23").Value = 6870
$ws.Range("R23").Value = 103009196900000

$ws.Range("A24").Value = "Casa"
$ws.Range("Q24").Value = 9273
$ws.Range("R24").Value = 0

$ws.Range("A25").Value = "Casa"
$ws.Range("O25").Value = "Leilão SFI Caixa"
$ws.Range("Q25").Value = 10845
$ws.Range("R25").Value = 40768

$ws.Range("A26").Value = "Casa"
$ws.Range("O26").Value = "Leilão SFI Caixa"
$ws.Range("Q26").Value = 8272
$ws.Range("R26").Value = 0

$ws.Range("A27").Value = "Apartamento"
$ws.Range("O27").Value = "Leilão SFI Caixa"
$ws.Range("Q27").Value = 23862
$ws.Range("R27").Value = 101000000100431

$ws.Range("A28").Value = "Casa"
$ws.Range("O28").Value = "Venda Direta"

$ws.Range("A29").Value = "Terreno"
$ws.Range("O29").Value = "Venda Direta"

$ws.Range("A30").Value = "Outros"
$ws.Range("O30").Value = "Venda Direta"

$ws.Range("A31").Value = "Terreno"
$ws.Range("O31").Value = "Venda Direta"
